# Auto-generated Excel COM-interop script
# Applies per-cell profit/price recalculation updates to the Leve tables
# across the ALC, ARM, BSM, CRP, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 51
$ws.Range("H51").Value = 7470.5
$ws.Range("I51").Value = 6900
$ws.Range("J51").Value = 7660.6665
$ws.Range("K51").Value = 6900
$ws.Range("L51").Value = 7660.6665
$ws.Range("M51").Value = -6416
$ws.Range("N51").Value = -8628.666499999999

# Row 62
$ws.Range("H62").Value = 6337.3076
$ws.Range("I62").Value = 7355.222
$ws.Range("J62").Value = 4047
$ws.Range("K62").Value = 7355.222
$ws.Range("L62").Value = 4047
$ws.Range("M62").Value = -6731.222
$ws.Range("N62").Value = -5295

# Row 65
$ws.Range("H65").Value = 6337.3076
$ws.Range("I65").Value = 7355.222
$ws.Range("J65").Value = 4047
$ws.Range("K65").Value = 36776.11
$ws.Range("L65").Value = 20235
$ws.Range("M65").Value = -33656.11
$ws.Range("N65").Value = -26475

# Row 74
$ws.Range("H74").Value = 14143.143
$ws.Range("I74").Value = 12003
$ws.Range("J74").Value = 14499.833
$ws.Range("K74").Value = 12003
$ws.Range("L74").Value = 14499.833
$ws.Range("M74").Value = -11067
$ws.Range("N74").Value = -16371.833

# Row 77
$ws.Range("H77").Value = 14143.143
$ws.Range("I77").Value = 12003
$ws.Range("J77").Value = 14499.833
$ws.Range("K77").Value = 60015
$ws.Range("L77").Value = 72499.16500000001
$ws.Range("M77").Value = -55335
$ws.Range("N77").Value = -81859.16500000001

# Row 137
$ws.Range("H137").Value = 2095.1365
$ws.Range("J137").Value = 3166
$ws.Range("L137").Value = 9498
$ws.Range("N137").Value = -14598


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 1775.7391
$ws.Range("I2").Value = 1230.5714
$ws.Range("J2").Value = 7500
$ws.Range("K2").Value = 1230.5714
$ws.Range("L2").Value = 7500
$ws.Range("M2").Value = -1117.5714
$ws.Range("N2").Value = -7726

# Row 32
$ws.Range("H32").Value = 7967.923
$ws.Range("I32").Value = 7967.923
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7967.923
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7680.923
$ws.Range("N32").ClearContents()

# Row 45
$ws.Range("H45").Value = 5767.5293
$ws.Range("I45").Value = 4860.5713
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 4860.5713
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -4483.5713
$ws.Range("N45").Value = -10754

# Row 61
$ws.Range("H61").Value = 22312.834
$ws.Range("I61").Value = 24775.4
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 24775.4
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -24563.4
$ws.Range("N61").Value = -10424

# Row 74
$ws.Range("H74").Value = 4310.7856
$ws.Range("J74").Value = 9377.799999999999
$ws.Range("L74").Value = 9377.799999999999
$ws.Range("N74").Value = -11125.8

# Row 77
$ws.Range("H77").Value = 4310.7856
$ws.Range("J77").Value = 9377.799999999999
$ws.Range("L77").Value = 46889
$ws.Range("N77").Value = -55625

# Row 116
$ws.Range("H116").Value = 1775.7391
$ws.Range("I116").Value = 1230.5714
$ws.Range("J116").Value = 7500
$ws.Range("K116").Value = 1230.5714
$ws.Range("L116").Value = 7500
$ws.Range("M116").Value = 1063.4286
$ws.Range("N116").Value = -12088

# Row 136
$ws.Range("H136").Value = 22312.834
$ws.Range("I136").Value = 24775.4
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 74326.20000000001
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -71776.20000000001
$ws.Range("N136").Value = -35100


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 1775.7391
$ws.Range("I3").Value = 1230.5714
$ws.Range("J3").Value = 7500
$ws.Range("K3").Value = 1230.5714
$ws.Range("L3").Value = 7500
$ws.Range("M3").Value = -1116.5714
$ws.Range("N3").Value = -7728

# Row 22
$ws.Range("H22").Value = 1387.25
$ws.Range("I22").Value = 1266.5
$ws.Range("K22").Value = 1266.5
$ws.Range("M22").Value = -1093.5

# Row 86
$ws.Range("H86").Value = 966.3333
$ws.Range("I86").Value = 999.5
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 999.5
$ws.Range("L86").Value = 900
$ws.Range("M86").Value = 123.5
$ws.Range("N86").Value = -3146

# Row 89
$ws.Range("H89").Value = 966.3333
$ws.Range("I89").Value = 999.5
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 4997.5
$ws.Range("L89").Value = 4500
$ws.Range("M89").Value = 618.5
$ws.Range("N89").Value = -15732

# Row 99
$ws.Range("H99").Value = 3989.5417
$ws.Range("I99").Value = 2963.6
$ws.Range("K99").Value = 2963.6
$ws.Range("M99").Value = -1465.6

# Row 134
$ws.Range("H134").Value = 4158.108
$ws.Range("I134").Value = 3848.5293
$ws.Range("K134").Value = 11545.5879
$ws.Range("M134").Value = -9010.5879


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 38467868
$ws.Range("I31").Value = 58826630
$ws.Range("K31").Value = 58826630
$ws.Range("M31").Value = -58826335

# Row 34
$ws.Range("H34").Value = 38467868
$ws.Range("I34").Value = 58826630
$ws.Range("K34").Value = 58826630
$ws.Range("M34").Value = -58826428

# Row 99
$ws.Range("H99").Value = 4743.769
$ws.Range("I99").Value = 4796
$ws.Range("K99").Value = 4796
$ws.Range("M99").Value = -3298

# Row 126
$ws.Range("H126").Value = 4743.769
$ws.Range("I126").Value = 4796
$ws.Range("K126").Value = 14388
$ws.Range("M126").Value = -11918

# Row 132
$ws.Range("H132").Value = 5899.8
$ws.Range("I132").Value = 3227
$ws.Range("K132").Value = 9681
$ws.Range("M132").Value = -7151

# Row 134
$ws.Range("H134").Value = 5046.7896
$ws.Range("I134").Value = 3406.3572
$ws.Range("K134").Value = 10219.0716
$ws.Range("M134").Value = -7684.071599999999


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 70
$ws.Range("H70").Value = 10674.2
$ws.Range("I70").Value = 7280.875
$ws.Range("K70").Value = 7280.875
$ws.Range("M70").Value = -7010.875

# Row 73
$ws.Range("H73").Value = 10674.2
$ws.Range("I73").Value = 7280.875
$ws.Range("K73").Value = 7280.875
$ws.Range("M73").Value = -6344.875

# Row 80
$ws.Range("H80").Value = 5870.923
$ws.Range("I80").Value = 5825.1665
$ws.Range("J80").Value = 5910.143
$ws.Range("K80").Value = 5825.1665
$ws.Range("L80").Value = 5910.143
$ws.Range("M80").Value = -4827.1665
$ws.Range("N80").Value = -7906.143

# Row 83
$ws.Range("H83").Value = 5870.923
$ws.Range("I83").Value = 5825.1665
$ws.Range("J83").Value = 5910.143
$ws.Range("K83").Value = 29125.8325
$ws.Range("L83").Value = 29550.715
$ws.Range("M83").Value = -24133.8325
$ws.Range("N83").Value = -39534.715

# Row 126
$ws.Range("H126").Value = 4181
$ws.Range("I126").Value = 3816
$ws.Range("K126").Value = 11448
$ws.Range("M126").Value = -8978


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 1442.5714
$ws.Range("I22").Value = 1424.625
$ws.Range("J22").Value = 1466.5
$ws.Range("K22").Value = 1424.625
$ws.Range("L22").Value = 1466.5
$ws.Range("M22").Value = -1129.625
$ws.Range("N22").Value = -2056.5

# Row 27
$ws.Range("H27").Value = 1442.5714
$ws.Range("I27").Value = 1424.625
$ws.Range("J27").Value = 1466.5
$ws.Range("K27").Value = 1424.625
$ws.Range("L27").Value = 1466.5
$ws.Range("M27").Value = -1317.625
$ws.Range("N27").Value = -1680.5

# Row 46
$ws.Range("H46").Value = 2320.2
$ws.Range("I46").Value = 2649.75
$ws.Range("J46").Value = 1002
$ws.Range("K46").Value = 2649.75
$ws.Range("L46").Value = 1002
$ws.Range("M46").Value = -2461.75
$ws.Range("N46").Value = -1378

# Row 108
$ws.Range("H108").Value = 74996
$ws.Range("J108").Value = 74996
$ws.Range("L108").Value = 74996
$ws.Range("N108").Value = -82676

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value = 4363.304
$ws.Range("I81").Value = 2548.75
$ws.Range("K81").Value = 5097.5
$ws.Range("M81").Value = -4036.5

# Row 84
$ws.Range("H84").Value = 4363.304
$ws.Range("I84").Value = 2548.75
$ws.Range("K84").Value = 25487.5
$ws.Range("M84").Value = -20183.5

# Row 122
$ws.Range("H122").Value = 4163.6875
$ws.Range("I122").Value = 3044.9565
$ws.Range("K122").Value = 9134.869499999999
$ws.Range("M122").Value = -6684.869499999999

